$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC1 (rows 10-13): step 2 (row 11) and the CAS/TJSeg alert swap places' content
$ws.Range("B11").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("D11").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"
$ws.Range("B12").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"

# TC2 (rows 20-23): step 2 and step 3 both now read "seleciona..."
$ws.Range("B21").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("B22").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"

# TC3 (rows 30-33): step 2 text and its expected result swap with the CAS alert
$ws.Range("B31").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("D31").Value = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar"
